$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.909.79'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '2.344.46'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''307.37'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('D6').Value = '''101.36'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('E7').Value = '  -4.55%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -2.21%  '
$ws.Range('D11').Value = '''52.68'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('E14').Value = '  -2.76%  '
$ws.Range('D15').Value = '2.711.93'
$ws.Range('E15').Value = '  +1.45%  '
$ws.Range('D16').Value = '''15.46'
$ws.Range('E16').Value = '  +3.13%  '
$ws.Range('D17').Value = '2.344.26'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').Value = '42.845.64'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').Value = '''6.26'
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('D21').Value = '''11.73'
$ws.Range('E21').Value = '  -6.51%  '
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('D23').Value = '''67.47'
$ws.Range('E23').Value = '  -1.44%  '
$ws.Range('D24').Value = '''237.42'
$ws.Range('E24').Value = '  -1.69%  '
$ws.Range('D25').Value = '''2.01'
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('E26').Value = '  -2.38%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '''25.34'
$ws.Range('E28').Value = '  +2.52%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '''35.42'
$ws.Range('E29').Value = '  -5.45%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '''9.41'
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '''2.08'
$ws.Range('E31').Value = '  -1.77%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '''160.03'
$ws.Range('E32').Value = '  -4.72%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''5.15'
$ws.Range('E34').Value = '  -3.34%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').Value = '''17.84'
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''2.48'
$ws.Range('E36').Value = '  +3.35%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.0729'
$ws.Range('E37').Value = '  -2.01%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '''3.02'
$ws.Range('E38').Value = '  -4.07%  '
$ws.Range('E39').Value = '  +7.25%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '''1.88'
$ws.Range('E40').Value = '  +1.75%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.104'
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '''0.113'
$ws.Range('E42').Value = '  -3.20%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').Value = '''2.51'
$ws.Range('E43').Value = '  +8.86%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.021.75'
$ws.Range('E44').Value = '  +2.40%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''19.08'
$ws.Range('E45').Value = '  -3.73%  '
$ws.Range('D46').Value = '''0.0285'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '''10.50'
$ws.Range('E47').Value = '  +6.90%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''3.01'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '''56.89'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').Value = '''2.92'
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.578.21'
$ws.Range('E51').Value = '  +1.47%  '
